# Auto-generated edit script applying the target diff to "Year 11.xlsx"-style workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 11.6
$ws.Range("H2").Value = 23.2
$ws.Range("I2").Value = 29
$ws.Range("J2").Value = 34.8
$ws.Range("K2").Value = 40.6
$ws.Range("L2").Value = 46.4
$ws.Range("M2").Value = 52.2
$ws.Range("N2").Value = 58
$ws.Range("O2").Value = 52.2
$ws.Range("P2").Value = 46.4
$ws.Range("Q2").Value = 40.6
$ws.Range("R2").Value = 29
$ws.Range("S2").Value = 17.4
$ws.Range("T2").Value = 11.6
$ws.Range("I3").Value = 23.2
$ws.Range("J3").Value = 34.8
$ws.Range("K3").Value = 46.4
$ws.Range("L3").Value = 52.2
$ws.Range("M3").Value = 58
$ws.Range("N3").Value = 46.4
$ws.Range("O3").Value = 40.6
$ws.Range("P3").Value = 29
$ws.Range("Q3").Value = 29
$ws.Range("R3").Value = 17.4
$ws.Range("S3").Value = 11.6
$ws.Range("J4").Value = 5.8
$ws.Range("K4").Value = 23.2
$ws.Range("L4").Value = 40.6
$ws.Range("M4").Value = 46.4
$ws.Range("N4").Value = 46.4
$ws.Range("O4").Value = 40.6
$ws.Range("P4").Value = 23.2
$ws.Range("Q4").Value = 11.6
$ws.Range("R4").Value = 5.8

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 62.3
$ws.Range("H2").Value = 10.2
$ws.Range("J2").Value = 111.934272013061
$ws.Range("K2").Value = 14.6
$ws.Range("L2").Value = 25.6
$ws.Range("M2").Value = 28.8
$ws.Range("N2").Value = 32
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 17.8
$ws.Range("Q2").Value = 14.6
$ws.Range("R2").Value = 43.9
$ws.Range("S2").Value = 57.8
$ws.Range("T2").Value = 43.6
$ws.Range("I3").Value = 23.2
$ws.Range("J3").Value = 34.8
$ws.Range("K3").Value = 46.4
$ws.Range("L3").Value = 52.2
$ws.Range("M3").Value = 34.6
$ws.Range("N3").Value = 20.4
$ws.Range("O3").Value = 105.8307927762476
$ws.Range("P3").Value = 0.4
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 17.4
$ws.Range("J4").Value = 15.8
$ws.Range("K4").Value = 23.2
$ws.Range("L4").Value = 40.6
$ws.Range("M4").Value = 28.18312417100182
$ws.Range("N4").Value = 46.4
$ws.Range("O4").Value = 40.6
$ws.Range("P4").Value = 23.2
$ws.Range("Q4").Value = 11.6
$ws.Range("R4").Value = 15.8

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("G2").Value = 181.677
$ws.Range("H2").Value = 191.775
$ws.Range("I2").Value = 191.775
$ws.Range("J2").Value = 302.5899292929304
$ws.Range("K2").Value = 317.0439292929304
$ws.Range("L2").Value = 342.3879292929304
$ws.Range("M2").Value = 370.8999292929304
$ws.Range("N2").Value = 402.5799292929304
$ws.Range("O2").Value = 423.3699292929305
$ws.Range("P2").Value = 440.9919292929305
$ws.Range("Q2").Value = 455.4459292929305
$ws.Range("R2").Value = 498.9069292929305
$ws.Range("S2").Value = 556.1289292929305
$ws.Range("I3").Value = 142.968
$ws.Range("J3").Value = 177.42
$ws.Range("K3").Value = 223.356
$ws.Range("L3").Value = 275.034
$ws.Range("M3").Value = 309.288
$ws.Range("N3").Value = 329.484
$ws.Range("O3").Value = 434.2564848484852
$ws.Range("P3").Value = 434.6524848484852
$ws.Range("Q3").Value = 437.6224848484852
$ws.Range("B4").Value = 159.3939393939379
$ws.Range("C4").Value = 139.6969696969689
$ws.Range("D4").Value = 139.6969696969689
$ws.Range("E4").Value = 139.6969696969689
$ws.Range("F4").Value = 139.6969696969689
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = 120
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 135.642
$ws.Range("K4").Value = 158.61
$ws.Range("L4").Value = 198.804
$ws.Range("M4").Value = 226.7052929292918
$ws.Range("N4").Value = 272.6412929292918
$ws.Range("O4").Value = 312.8352929292918
$ws.Range("P4").Value = 335.8032929292918
$ws.Range("Q4").Value = 347.2872929292918
$ws.Range("R4").Value = 362.9292929292918
$ws.Range("S4").Value = 362.9292929292918
$ws.Range("T4").Value = 231.6161616161605
$ws.Range("U4").Value = 231.6161616161605
$ws.Range("V4").Value = 231.6161616161605
$ws.Range("W4").Value = 231.6161616161605
$ws.Range("X4").Value = 231.6161616161605
$ws.Range("Y4").Value = 192.2222222222222

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("G2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("O3").Value = 65.23079277624761
$ws.Range("K4").Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 2.2
$ws.Range("M2").Value = 0
$ws.Range("R2").Value = 48.7
$ws.Range("I3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("Q4").Value = 0

$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Range("J2").Value = 57.63427201306106
$ws.Range("Q2").Value = 0
$ws.Range("S2").Value = 53.6
$ws.Range("M3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("S3").Value = 11.6
$ws.Range("M4").Value = 5.183124171001815
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 77388.6679767371
$ws.Range("C2").Value = 2000
$ws.Range("D2").Value = 10661.60663177227
$ws.Range("E2").Value = 1620
$ws.Range("F2").Value = 20022.81757366193

$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 58
$ws.Range("B4").Value = 10
$ws.Range("D4").Value = 10

